# Apply cell content updates per the cryptos list refresh (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.500.85"
$ws.Range("E2").Value = "  +2.85%  "
$ws.Range("D3").Value = "3.022.52"
$ws.Range("E3").Value = "  +5.70%  "
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "508.04"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.22"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +6.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.428"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +4.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.14"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.106"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.365"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +7.75%  "
$ws.Range("D12").Value = "3.542.97"
$ws.Range("E12").Value = "  +5.86%  "
$ws.Range("E13").Value = "  +3.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.81"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000160"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +5.32%  "
$ws.Range("D16").Value = "56.569.02"
$ws.Range("E16").Value = "  +2.75%  "
$ws.Range("D17").Value = "3.024.96"
$ws.Range("E17").Value = "  +5.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.87"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.98"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +7.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.02"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +7.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "331.08"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +8.12%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.497"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +6.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.36"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +6.47%  "
$ws.Range("D25").Value = "3.155.77"
$ws.Range("E25").Value = "  +6.34%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.165"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.67%  "
$ws.Range("D28").Value = "0.0₃0913"
$ws.Range("E28").Value = "  +14.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.27"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.80"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.78"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +5.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.49"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +6.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.14"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +5.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "153.13"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +5.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.45"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.63"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +13.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.78"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +6.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.21"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0664"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.38%  "
$ws.Range("D40").Value = "3.054.83"
$ws.Range("E40").Value = "  +5.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.62"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.34%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.657"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.99%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.78"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +6.80%  "
$ws.Range("D45").Value = "2.213.18"
$ws.Range("E45").Value = "  +7.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0248"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +11.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.34"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.928"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.37%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.80"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.22%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.51"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +7.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0861"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +4.95%  "
